# Update marksheet "Marking" and "Total" rows (quiz sheet) to reflect
# corrected per-question mark and recalculated total score.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row (row 11): points per correct answer 3 -> 5
$ws.Range("B11").Value = 5

# Total row (row 12): total score 75 -> 125, shown as "125/140"
$ws.Range("B12").Value = 125
$ws.Range("E12").Value = "125/140"
